$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formats for the new row's cells from the closest matching existing
# --- rows so the new cells pick up the same visual style (wrap text font etc).
# C18 should look like the other "Access URL" cells that use the small
# grey Arial font (style used on C8 / C13).
$ws.Cells.Item(8, 3).Copy()
$ws.Cells.Item(18, 3).PasteSpecial(-4122)

# F18 / G18 should look like the wrapped JSON cells used throughout the
# sheet (style used on F17 / G17).
$ws.Cells.Item(17, 6).Copy()
$ws.Cells.Item(18, 6).PasteSpecial(-4122)

$ws.Cells.Item(17, 7).Copy()
$ws.Cells.Item(18, 7).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Request body / response JSON payloads (kept as here-strings so the
# --- literal text, quotes, ampersand and newlines are preserved exactly).
$postBody = @'
{
 "userUid":"5U9wbAfe66RrFRibw07u9qBkDQJ2",
 "postTitle":"Post 1",
 "postImage":"https://firebasestorage.googleapis.com/v0/b/discussion-manager.appspot.com/o/annie-spratt-QckxruozjRg-unsplash.jpg?alt=media&token=922ba71c-45dd-4f46-85ba-6030eb80cea9",
 "sentTimeDate":"7:50",
 "postDesc":"New post desc"
}

'@

$postResponse = @'
{
    "responseMessage": "Uploaded the post",
    "responseCode": 2
}
'@

# --- Fill in the new "Posts" end point row (row 18).
$ws.Cells.Item(18, 2).Value = "Posts"
$ws.Cells.Item(18, 3).Value = "http://localhost:8000/postsServices/addPost"
$ws.Cells.Item(18, 4).Value = "to upload the posts"
$ws.Cells.Item(18, 5).Value = "POST"
$ws.Cells.Item(18, 6).Value = $postBody
$ws.Cells.Item(18, 7).Value = $postResponse

# Row height matches the other wrapped-text rows on the sheet.
$ws.Rows.Item(18).RowHeight = 195

# --- Update the view so the newly added row is visible / selected, same as
# --- what the author would see right after typing the new row in.
$ws.Range("F30").Select()
